# assigneng.xlsx - "modified and updated the insertexcel and form2"
#
# The sheet used to hold 7 data rows (rows 2-8). The edit trims it down to a
# single data row (row 2) and rewrites that row's values - including turning
# the "Due Date" cell from a real date value into a literal text string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 3-8 entirely - only the header (row 1) and one data row (row 2)
# remain, shrinking the used range down to A1:G2.
$ws.Rows("3:8").Delete()

# Row 2 stays in place for AssignId (column A); the rest of the row gets new
# values.
$ws.Range("B2").Value = "sekar"
$ws.Range("C2").Value = "jjjjjj"

# D2 used to be a real date (serial number, formatted as a date). Force it to
# be stored as literal text "2018-05-12" instead of letting Excel parse it
# back into a date serial, then drop the temporary formatting so the cell
# ends up with no explicit style, same as the rest of the row.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2018-05-12"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "Once a Day"
$ws.Range("F2").Value = "Tender"
$ws.Range("G2").Value = "ds"

# Match the author's final selection in the sheet.
$ws.Range("C13").Select()
